$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - these cells store numeric-looking text, so force
# the Text number format first to keep Excel from auto-converting the
# assigned string into a numeric value (preserves the original inlineStr type).
$priceUpdates = @{
    "D3"  = "23.06"
    "D4"  = "5.408"
    "D6"  = "3.390"
    "D7"  = "0.8098"
    "D8"  = "0.9251"
    "D10" = "0.07427"
    "D11" = "0.03370"
    "D12" = "0.03032"
    "D13" = "0.09358"
    "D14" = "3.969"
    "D15" = "0.001601"
    "D16" = "0.04826"
    "D18" = "0.005301"
    "D19" = "0.004154"
    "D20" = "0.0009814"
    "D22" = "3.655"
    "D23" = "6.441"
    "D24" = "2.186"
    "D26" = "0.1294"
    "D40" = "0.03969"
    "D41" = "0.006462"
    "D42" = "0.1073"
    "D43" = "0.002901"
    "D44" = "0.006714"
    "D45" = "0.00005199"
    "D49" = "0.002265"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Column E (Volume(1h)) updates - plain text, no numeric coercion needed.
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
